$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '68.662.63'
$ws.Cells.Item(2, 5).Value = '  +2.32%  '

$ws.Cells.Item(3, 4).Value = '3.749.31'
$ws.Cells.Item(3, 5).Value = '  +1.88%  '

$ws.Cells.Item(4, 5).Value = '  +0.13%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '601.06'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.77%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '168.70'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +1.26%  '

$ws.Cells.Item(7, 4).Value = '3.748.64'
$ws.Cells.Item(7, 5).Value = '  +1.86%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.535'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  +3.08%  '

$ws.Cells.Item(10, 5).Value = '  +3.80%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '6.32'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +3.14%  '

$ws.Cells.Item(12, 5).Value = '  +0.94%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '38.27'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.28%  '

$ws.Cells.Item(14, 5).Value = '  +3.48%  '

$ws.Cells.Item(15, 4).Value = '4.377.28'
$ws.Cells.Item(15, 5).Value = '  +2.03%  '

$ws.Cells.Item(16, 4).Value = '3.754.77'
$ws.Cells.Item(16, 5).Value = '  +2.25%  '

$ws.Cells.Item(17, 4).Value = '68.696.20'
$ws.Cells.Item(17, 5).Value = '  +2.39%  '

$ws.Cells.Item(18, 5).Value = '  +3.40%  '

$ws.Cells.Item(19, 5).Value = '  +0.93%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '17.11'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +1.70%  '

$ws.Cells.Item(21, 5).Value = '  +20.61%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '495.24'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +2.63%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.728'
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +1.96%  '

$ws.Cells.Item(24, 5).Value = '  +9.63%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '85.25'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +0.91%  '

$ws.Cells.Item(26, 5).Value = '  +1.64%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '12.37'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +2.52%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '10.25'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +3.41%  '

$ws.Cells.Item(29, 5).Value = '  +0.43%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '2.53'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +8.29%  '

$ws.Cells.Item(31, 5).Value = '  +2.90%  '

$ws.Cells.Item(32, 5).Value = '  +3.17%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '31.86'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +0.60%  '

$ws.Cells.Item(34, 4).Value = '3.895.06'
$ws.Cells.Item(34, 5).Value = '  +2.15%  '

$ws.Cells.Item(35, 4).Value = '3.684.78'
$ws.Cells.Item(35, 5).Value = '  +1.89%  '

$ws.Cells.Item(36, 5).Value = '  +2.55%  '

$ws.Cells.Item(37, 5).Value = '  +0.19%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '1.01'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +2.78%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.85'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +2.42%  '

$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.133'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.89%  '

$ws.Cells.Item(41, 5).Value = '  +1.54%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '441.35'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +1.38%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '48.75'
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.40%  '

$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.90'
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +5.51%  '

$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '1.98'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  +3.10%  '

$ws.Cells.Item(46, 5).Value = '  +2.82%  '

$ws.Cells.Item(47, 5).Value = '  +0.03%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '40.33'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +2.64%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '141.52'
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +0.11%  '

$ws.Cells.Item(50, 4).Value = '2.794.99'
$ws.Cells.Item(50, 5).Value = '  +1.68%  '

$ws.Cells.Item(51, 5).Value = '  +3.49%  '
